$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.842447280883789
$ws.Range("B1").Value = 5.811890125274658
$ws.Range("C1").Value = 4.773184776306152
$ws.Range("D1").Value = 5.346162796020508
$ws.Range("E1").Value = 5.633428573608398
